# "changes during the course"
#
# Reposition/resize four custom shapes (shape ids 195, 196, 197, 198) that
# sit on slide 8 (sldId 263) of the presentation -- part of a "Geom-specific"
# diagram made of two arrow/line shapes and two small rectangle labels.
#
# Target off/ext values (EMU), taken from the authored OOXML:
#   id195: off (3444480,1587600) -> (2875680,1587600); ext (3732120,503640) -> (4300920,503640)
#   id196: off (2459160,1437120) -> (1891080,1441800); ext unchanged
#   id197: off (3249360,974160)  -> (3154360,974160);  ext unchanged
#   id198: off (2277000,1143360) -> (2182000,1143360); ext unchanged
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points (1 pt =
# 12700 EMU) and are stored as single-precision floats, so the literals below
# are chosen to round-trip to the exact target EMU values after the
# points -> float32 -> EMU conversion the host performs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

foreach ($sh in $s.Shapes) {
    switch ($sh.Id) {
        195 {
            $sh.Left  = 226.43149606299212
            $sh.Width = 338.6551061102362
        }
        196 {
            $sh.Left = 148.90393900787402
            $sh.Top  = 113.52756205511811
        }
        197 {
            $sh.Left = 248.3748101496063
        }
        198 {
            $sh.Left = 171.81102762204722
        }
    }
}
